$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C ("Quantité Ingrédients") to hold the
# numeric "Quantite" values.
$ws.Columns("C:C").Insert()

# The inserted column should keep the same width as column B.
$ws.Columns("C:C").ColumnWidth = 36

# Header for the new column.
$ws.Range("C1").Value = "Quantite"

# Numeric quantity values for each ingredient row (rows 2-10).
$ws.Range("C2").Value = 2
$ws.Range("C3").Value = 0
$ws.Range("C4").Value = 125
$ws.Range("C5").Value = 1
$ws.Range("C6").Value = 1
$ws.Range("C7").Value = 0.5
$ws.Range("C8").Value = 0
$ws.Range("C9").Value = 0.5
$ws.Range("C10").Value = 150

# Update the view: move the active selection to C10.
[void]$ws.Range("C10").Select()
